$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Converted Data")

# --- W7 / AA7 updates ---
$ws.Range("W7").Value = 0
$ws.Range("AA7").Value = 12

# --- AA column recomputed ratios (denominator changed 26->24 days etc.) ---
$ws.Range("AA25:AA37").Value = 0.1833333333333333
$ws.Range("AA38:AA72").Value = 0.8333333333249998
$ws.Range("AA73").Value = 0.7499999999916667
$ws.Range("AA74:AA87").Value = 0.6249999999916667
$ws.Range("AA88:AA91").Value = 0.5000000000000001
$ws.Range("AA92:AA108").Value = 0.3083333333333333
$ws.Range("AA109:AA132").Value = 0.1666666666666667
$ws.Range("AA133:AA221").Value = 0.25

# --- New rows 222-233: dates 9/30/2020 .. 10/11/2020 ---
$a = $ws.Cells.Item(222, 1)
$a.NumberFormat = "@"
$a.Value = "9/30/2020"
$ws.Cells.Item(222, 2).Value = 0
$ws.Cells.Item(222, 3).Value = 0
$ws.Cells.Item(222, 4).Value = 1
$ws.Cells.Item(222, 5).Value = 0
$ws.Cells.Item(222, 6).Value = 1
$ws.Cells.Item(222, 7).Value = 0
$ws.Cells.Item(222, 8).Value = 0
$ws.Cells.Item(222, 9).Value = 0
$ws.Cells.Item(222, 10).Value = 0
$ws.Cells.Item(222, 11).Value = 0
$ws.Cells.Item(222, 12).Value = 0
$ws.Cells.Item(222, 13).Value = 0
$ws.Cells.Item(222, 14).Value = 0
$ws.Cells.Item(222, 15).Value = 0
$ws.Cells.Item(222, 16).Value = 0
$ws.Cells.Item(222, 17).Value = 0
$ws.Cells.Item(222, 18).Value = 0
$ws.Cells.Item(222, 19).Value = 0
$ws.Cells.Item(222, 20).Value = 0
$ws.Cells.Item(222, 21).Value = 0
$ws.Cells.Item(222, 22).Value = 0
$ws.Cells.Item(222, 23).Value = 0
$ws.Cells.Item(222, 24).Value = 0
$ws.Cells.Item(222, 25).Value = 1
$ws.Cells.Item(222, 26).Value = 1
$ws.Cells.Item(222, 27).Value = 0.25

$a = $ws.Cells.Item(223, 1)
$a.NumberFormat = "@"
$a.Value = "10/1/2020"
$ws.Cells.Item(223, 2).Value = 0
$ws.Cells.Item(223, 3).Value = 0
$ws.Cells.Item(223, 4).Value = 1
$ws.Cells.Item(223, 5).Value = 0
$ws.Cells.Item(223, 6).Value = 1
$ws.Cells.Item(223, 7).Value = 0
$ws.Cells.Item(223, 8).Value = 0
$ws.Cells.Item(223, 9).Value = 0
$ws.Cells.Item(223, 10).Value = 0
$ws.Cells.Item(223, 11).Value = 0
$ws.Cells.Item(223, 12).Value = 0
$ws.Cells.Item(223, 13).Value = 0
$ws.Cells.Item(223, 14).Value = 0
$ws.Cells.Item(223, 15).Value = 0
$ws.Cells.Item(223, 16).Value = 0
$ws.Cells.Item(223, 17).Value = 0
$ws.Cells.Item(223, 18).Value = 0
$ws.Cells.Item(223, 19).Value = 0
$ws.Cells.Item(223, 20).Value = 0
$ws.Cells.Item(223, 21).Value = 0
$ws.Cells.Item(223, 22).Value = 0
$ws.Cells.Item(223, 23).Value = 0
$ws.Cells.Item(223, 24).Value = 0
$ws.Cells.Item(223, 25).Value = 1
$ws.Cells.Item(223, 26).Value = 1
$ws.Cells.Item(223, 27).Value = 0.25

$a = $ws.Cells.Item(224, 1)
$a.NumberFormat = "@"
$a.Value = "10/2/2020"
$ws.Cells.Item(224, 2).Value = 0
$ws.Cells.Item(224, 3).Value = 0
$ws.Cells.Item(224, 4).Value = 1
$ws.Cells.Item(224, 5).Value = 0
$ws.Cells.Item(224, 6).Value = 1
$ws.Cells.Item(224, 7).Value = 0
$ws.Cells.Item(224, 8).Value = 0
$ws.Cells.Item(224, 9).Value = 0
$ws.Cells.Item(224, 10).Value = 0
$ws.Cells.Item(224, 11).Value = 0
$ws.Cells.Item(224, 12).Value = 0
$ws.Cells.Item(224, 13).Value = 0
$ws.Cells.Item(224, 14).Value = 0
$ws.Cells.Item(224, 15).Value = 0
$ws.Cells.Item(224, 16).Value = 0
$ws.Cells.Item(224, 17).Value = 0
$ws.Cells.Item(224, 18).Value = 0
$ws.Cells.Item(224, 19).Value = 0
$ws.Cells.Item(224, 20).Value = 0
$ws.Cells.Item(224, 21).Value = 0
$ws.Cells.Item(224, 22).Value = 0
$ws.Cells.Item(224, 23).Value = 0
$ws.Cells.Item(224, 24).Value = 0
$ws.Cells.Item(224, 25).Value = 1
$ws.Cells.Item(224, 26).Value = 1
$ws.Cells.Item(224, 27).Value = 0.25

$a = $ws.Cells.Item(225, 1)
$a.NumberFormat = "@"
$a.Value = "10/3/2020"
$ws.Cells.Item(225, 2).Value = 0
$ws.Cells.Item(225, 3).Value = 0
$ws.Cells.Item(225, 4).Value = 1
$ws.Cells.Item(225, 5).Value = 0
$ws.Cells.Item(225, 6).Value = 1
$ws.Cells.Item(225, 7).Value = 0
$ws.Cells.Item(225, 8).Value = 0
$ws.Cells.Item(225, 9).Value = 0
$ws.Cells.Item(225, 10).Value = 0
$ws.Cells.Item(225, 11).Value = 0
$ws.Cells.Item(225, 12).Value = 0
$ws.Cells.Item(225, 13).Value = 0
$ws.Cells.Item(225, 14).Value = 0
$ws.Cells.Item(225, 15).Value = 0
$ws.Cells.Item(225, 16).Value = 0
$ws.Cells.Item(225, 17).Value = 0
$ws.Cells.Item(225, 18).Value = 0
$ws.Cells.Item(225, 19).Value = 0
$ws.Cells.Item(225, 20).Value = 0
$ws.Cells.Item(225, 21).Value = 0
$ws.Cells.Item(225, 22).Value = 0
$ws.Cells.Item(225, 23).Value = 0
$ws.Cells.Item(225, 24).Value = 0
$ws.Cells.Item(225, 25).Value = 1
$ws.Cells.Item(225, 26).Value = 1
$ws.Cells.Item(225, 27).Value = 0.25

$a = $ws.Cells.Item(226, 1)
$a.NumberFormat = "@"
$a.Value = "10/4/2020"
$ws.Cells.Item(226, 2).Value = 0
$ws.Cells.Item(226, 3).Value = 0
$ws.Cells.Item(226, 4).Value = 1
$ws.Cells.Item(226, 5).Value = 0
$ws.Cells.Item(226, 6).Value = 1
$ws.Cells.Item(226, 7).Value = 0
$ws.Cells.Item(226, 8).Value = 0
$ws.Cells.Item(226, 9).Value = 0
$ws.Cells.Item(226, 10).Value = 0
$ws.Cells.Item(226, 11).Value = 0
$ws.Cells.Item(226, 12).Value = 0
$ws.Cells.Item(226, 13).Value = 0
$ws.Cells.Item(226, 14).Value = 0
$ws.Cells.Item(226, 15).Value = 0
$ws.Cells.Item(226, 16).Value = 0
$ws.Cells.Item(226, 17).Value = 0
$ws.Cells.Item(226, 18).Value = 0
$ws.Cells.Item(226, 19).Value = 0
$ws.Cells.Item(226, 20).Value = 0
$ws.Cells.Item(226, 21).Value = 0
$ws.Cells.Item(226, 22).Value = 0
$ws.Cells.Item(226, 23).Value = 0
$ws.Cells.Item(226, 24).Value = 0
$ws.Cells.Item(226, 25).Value = 1
$ws.Cells.Item(226, 26).Value = 1
$ws.Cells.Item(226, 27).Value = 0.25

$a = $ws.Cells.Item(227, 1)
$a.NumberFormat = "@"
$a.Value = "10/5/2020"
$ws.Cells.Item(227, 2).Value = 0
$ws.Cells.Item(227, 3).Value = 0
$ws.Cells.Item(227, 4).Value = 1
$ws.Cells.Item(227, 5).Value = 0
$ws.Cells.Item(227, 6).Value = 1
$ws.Cells.Item(227, 7).Value = 0
$ws.Cells.Item(227, 8).Value = 0
$ws.Cells.Item(227, 9).Value = 0
$ws.Cells.Item(227, 10).Value = 0
$ws.Cells.Item(227, 11).Value = 0
$ws.Cells.Item(227, 12).Value = 0
$ws.Cells.Item(227, 13).Value = 0
$ws.Cells.Item(227, 14).Value = 0
$ws.Cells.Item(227, 15).Value = 0
$ws.Cells.Item(227, 16).Value = 0
$ws.Cells.Item(227, 17).Value = 0
$ws.Cells.Item(227, 18).Value = 0
$ws.Cells.Item(227, 19).Value = 0
$ws.Cells.Item(227, 20).Value = 0
$ws.Cells.Item(227, 21).Value = 0
$ws.Cells.Item(227, 22).Value = 0
$ws.Cells.Item(227, 23).Value = 0
$ws.Cells.Item(227, 24).Value = 0
$ws.Cells.Item(227, 25).Value = 1
$ws.Cells.Item(227, 26).Value = 1
$ws.Cells.Item(227, 27).Value = 0.25

$a = $ws.Cells.Item(228, 1)
$a.NumberFormat = "@"
$a.Value = "10/6/2020"
$ws.Cells.Item(228, 2).Value = 0
$ws.Cells.Item(228, 3).Value = 0
$ws.Cells.Item(228, 4).Value = 1
$ws.Cells.Item(228, 5).Value = 0
$ws.Cells.Item(228, 6).Value = 1
$ws.Cells.Item(228, 7).Value = 0
$ws.Cells.Item(228, 8).Value = 0
$ws.Cells.Item(228, 9).Value = 0
$ws.Cells.Item(228, 10).Value = 0
$ws.Cells.Item(228, 11).Value = 0
$ws.Cells.Item(228, 12).Value = 0
$ws.Cells.Item(228, 13).Value = 0
$ws.Cells.Item(228, 14).Value = 0
$ws.Cells.Item(228, 15).Value = 0
$ws.Cells.Item(228, 16).Value = 0
$ws.Cells.Item(228, 17).Value = 0
$ws.Cells.Item(228, 18).Value = 0
$ws.Cells.Item(228, 19).Value = 0
$ws.Cells.Item(228, 20).Value = 0
$ws.Cells.Item(228, 21).Value = 0
$ws.Cells.Item(228, 22).Value = 0
$ws.Cells.Item(228, 23).Value = 0
$ws.Cells.Item(228, 24).Value = 0
$ws.Cells.Item(228, 25).Value = 1
$ws.Cells.Item(228, 26).Value = 1
$ws.Cells.Item(228, 27).Value = 0.25

$a = $ws.Cells.Item(229, 1)
$a.NumberFormat = "@"
$a.Value = "10/7/2020"
$ws.Cells.Item(229, 2).Value = 0
$ws.Cells.Item(229, 3).Value = 0
$ws.Cells.Item(229, 4).Value = 1
$ws.Cells.Item(229, 5).Value = 0
$ws.Cells.Item(229, 6).Value = 1
$ws.Cells.Item(229, 7).Value = 0
$ws.Cells.Item(229, 8).Value = 0
$ws.Cells.Item(229, 9).Value = 0
$ws.Cells.Item(229, 10).Value = 0
$ws.Cells.Item(229, 11).Value = 0
$ws.Cells.Item(229, 12).Value = 0
$ws.Cells.Item(229, 13).Value = 0
$ws.Cells.Item(229, 14).Value = 0
$ws.Cells.Item(229, 15).Value = 0
$ws.Cells.Item(229, 16).Value = 0
$ws.Cells.Item(229, 17).Value = 0
$ws.Cells.Item(229, 18).Value = 0
$ws.Cells.Item(229, 19).Value = 0
$ws.Cells.Item(229, 20).Value = 0
$ws.Cells.Item(229, 21).Value = 0
$ws.Cells.Item(229, 22).Value = 0
$ws.Cells.Item(229, 23).Value = 0
$ws.Cells.Item(229, 24).Value = 0
$ws.Cells.Item(229, 25).Value = 1
$ws.Cells.Item(229, 26).Value = 1
$ws.Cells.Item(229, 27).Value = 0.25

$a = $ws.Cells.Item(230, 1)
$a.NumberFormat = "@"
$a.Value = "10/8/2020"
$ws.Cells.Item(230, 2).Value = 0
$ws.Cells.Item(230, 3).Value = 0
$ws.Cells.Item(230, 4).Value = 1
$ws.Cells.Item(230, 5).Value = 0
$ws.Cells.Item(230, 6).Value = 1
$ws.Cells.Item(230, 7).Value = 0
$ws.Cells.Item(230, 8).Value = 0
$ws.Cells.Item(230, 9).Value = 0
$ws.Cells.Item(230, 10).Value = 0
$ws.Cells.Item(230, 11).Value = 0
$ws.Cells.Item(230, 12).Value = 0
$ws.Cells.Item(230, 13).Value = 0
$ws.Cells.Item(230, 14).Value = 0
$ws.Cells.Item(230, 15).Value = 0
$ws.Cells.Item(230, 16).Value = 0
$ws.Cells.Item(230, 17).Value = 0
$ws.Cells.Item(230, 18).Value = 0
$ws.Cells.Item(230, 19).Value = 0
$ws.Cells.Item(230, 20).Value = 0
$ws.Cells.Item(230, 21).Value = 0
$ws.Cells.Item(230, 22).Value = 0
$ws.Cells.Item(230, 23).Value = 0
$ws.Cells.Item(230, 24).Value = 0
$ws.Cells.Item(230, 25).Value = 1
$ws.Cells.Item(230, 26).Value = 1
$ws.Cells.Item(230, 27).Value = 0.25

$a = $ws.Cells.Item(231, 1)
$a.NumberFormat = "@"
$a.Value = "10/9/2020"
$ws.Cells.Item(231, 2).Value = 0
$ws.Cells.Item(231, 3).Value = 0
$ws.Cells.Item(231, 4).Value = 1
$ws.Cells.Item(231, 5).Value = 0
$ws.Cells.Item(231, 6).Value = 1
$ws.Cells.Item(231, 7).Value = 0
$ws.Cells.Item(231, 8).Value = 0
$ws.Cells.Item(231, 9).Value = 0
$ws.Cells.Item(231, 10).Value = 0
$ws.Cells.Item(231, 11).Value = 0
$ws.Cells.Item(231, 12).Value = 0
$ws.Cells.Item(231, 13).Value = 0
$ws.Cells.Item(231, 14).Value = 0
$ws.Cells.Item(231, 15).Value = 0
$ws.Cells.Item(231, 16).Value = 0
$ws.Cells.Item(231, 17).Value = 0
$ws.Cells.Item(231, 18).Value = 0
$ws.Cells.Item(231, 19).Value = 0
$ws.Cells.Item(231, 20).Value = 0
$ws.Cells.Item(231, 21).Value = 0
$ws.Cells.Item(231, 22).Value = 0
$ws.Cells.Item(231, 23).Value = 0
$ws.Cells.Item(231, 24).Value = 0
$ws.Cells.Item(231, 25).Value = 1
$ws.Cells.Item(231, 26).Value = 1
$ws.Cells.Item(231, 27).Value = 0.25

$a = $ws.Cells.Item(232, 1)
$a.NumberFormat = "@"
$a.Value = "10/10/2020"
$ws.Cells.Item(232, 2).Value = 0
$ws.Cells.Item(232, 3).Value = 0
$ws.Cells.Item(232, 4).Value = 1
$ws.Cells.Item(232, 5).Value = 0
$ws.Cells.Item(232, 6).Value = 1
$ws.Cells.Item(232, 7).Value = 0
$ws.Cells.Item(232, 8).Value = 0
$ws.Cells.Item(232, 9).Value = 0
$ws.Cells.Item(232, 10).Value = 0
$ws.Cells.Item(232, 11).Value = 0
$ws.Cells.Item(232, 12).Value = 0
$ws.Cells.Item(232, 13).Value = 0
$ws.Cells.Item(232, 14).Value = 0
$ws.Cells.Item(232, 15).Value = 0
$ws.Cells.Item(232, 16).Value = 0
$ws.Cells.Item(232, 17).Value = 0
$ws.Cells.Item(232, 18).Value = 0
$ws.Cells.Item(232, 19).Value = 0
$ws.Cells.Item(232, 20).Value = 0
$ws.Cells.Item(232, 21).Value = 0
$ws.Cells.Item(232, 22).Value = 0
$ws.Cells.Item(232, 23).Value = 0
$ws.Cells.Item(232, 24).Value = 0
$ws.Cells.Item(232, 25).Value = 1
$ws.Cells.Item(232, 26).Value = 1
$ws.Cells.Item(232, 27).Value = 0.25

$a = $ws.Cells.Item(233, 1)
$a.NumberFormat = "@"
$a.Value = "10/11/2020"
$ws.Cells.Item(233, 2).Value = 0
$ws.Cells.Item(233, 3).Value = 0
$ws.Cells.Item(233, 4).Value = 1
$ws.Cells.Item(233, 5).Value = 0
$ws.Cells.Item(233, 6).Value = 1
$ws.Cells.Item(233, 7).Value = 0
$ws.Cells.Item(233, 8).Value = 0
$ws.Cells.Item(233, 9).Value = 0
$ws.Cells.Item(233, 10).Value = 0
$ws.Cells.Item(233, 11).Value = 0
$ws.Cells.Item(233, 12).Value = 0
$ws.Cells.Item(233, 13).Value = 0
$ws.Cells.Item(233, 14).Value = 0
$ws.Cells.Item(233, 15).Value = 0
$ws.Cells.Item(233, 16).Value = 0
$ws.Cells.Item(233, 17).Value = 0
$ws.Cells.Item(233, 18).Value = 0
$ws.Cells.Item(233, 19).Value = 0
$ws.Cells.Item(233, 20).Value = 0
$ws.Cells.Item(233, 21).Value = 0
$ws.Cells.Item(233, 22).Value = 0
$ws.Cells.Item(233, 23).Value = 0
$ws.Cells.Item(233, 24).Value = 0
$ws.Cells.Item(233, 25).Value = 1
$ws.Cells.Item(233, 26).Value = 1
$ws.Cells.Item(233, 27).Value = 0.25
